# Natmi following Dr Hou advice
#
# The Sema3d -> Nrp1 ligand-receptor pair table gains a third cluster, "ECs"
# (endothelial cells), alongside the existing "FAPs" and "sCs" clusters, both as
# a sending cluster and as a target cluster. The 2x3 combinations become 3x3,
# so rows 2-7 are refreshed with the recomputed statistics and three new rows
# (8-10) are appended for the "sCs" sending-cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.03373200000000001
$ws.Cells.Item(2, 8).Value = 0.101196
$ws.Cells.Item(2, 9).Value = 0.002572411152219347
$ws.Cells.Item(2, 10).Value = 0.002572411152219347
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 102.8289443333334
$ws.Cells.Item(2, 14).Value = 308.486833
$ws.Cells.Item(2, 15).Value = 0.5559120396302444
$ws.Cells.Item(2, 16).Value = 0.5559120396302443
$ws.Cells.Item(2, 17).Value = 3.468625950252001
$ws.Cells.Item(2, 18).Value = 31.21763355226801
$ws.Cells.Item(2, 19).Value = 0.001430034330397844
$ws.Cells.Item(2, 20).Value = 0.001430034330397844

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.03373200000000001
$ws.Cells.Item(3, 8).Value = 0.101196
$ws.Cells.Item(3, 9).Value = 0.002572411152219347
$ws.Cells.Item(3, 10).Value = 0.002572411152219347
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 63.66262833333334
$ws.Cells.Item(3, 14).Value = 190.987885
$ws.Cells.Item(3, 15).Value = 0.3441717873742006
$ws.Cells.Item(3, 16).Value = 0.3441717873742006
$ws.Cells.Item(3, 17).Value = 2.14746777894
$ws.Cells.Item(3, 18).Value = 19.32721001046
$ws.Cells.Item(3, 19).Value = 0.0008853513441206594
$ws.Cells.Item(3, 20).Value = 0.0008853513441206594

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.03373200000000001
$ws.Cells.Item(4, 8).Value = 0.101196
$ws.Cells.Item(4, 9).Value = 0.002572411152219347
$ws.Cells.Item(4, 10).Value = 0.002572411152219347
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 18.481835
$ws.Cells.Item(4, 14).Value = 55.445505
$ws.Cells.Item(4, 15).Value = 0.09991617299555507
$ws.Cells.Item(4, 16).Value = 0.09991617299555505
$ws.Cells.Item(4, 17).Value = 0.6234292582200001
$ws.Cells.Item(4, 18).Value = 5.61086332398
$ws.Cells.Item(4, 19).Value = 0.0002570254777008434
$ws.Cells.Item(4, 20).Value = 0.0002570254777008433

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.484157
$ws.Cells.Item(5, 8).Value = 37.452471
$ws.Cells.Item(5, 9).Value = 0.952045081609665
$ws.Cells.Item(5, 10).Value = 0.9520450816096651
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 102.8289443333334
$ws.Cells.Item(5, 14).Value = 308.486833
$ws.Cells.Item(5, 15).Value = 0.5559120396302444
$ws.Cells.Item(5, 16).Value = 0.5559120396302443
$ws.Cells.Item(5, 17).Value = 1283.732685201594
$ws.Cells.Item(5, 18).Value = 11553.59416681435
$ws.Cells.Item(5, 19).Value = 0.5292533231375713
$ws.Cells.Item(5, 20).Value = 0.5292533231375712

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.484157
$ws.Cells.Item(6, 8).Value = 37.452471
$ws.Cells.Item(6, 9).Value = 0.952045081609665
$ws.Cells.Item(6, 10).Value = 0.9520450816096651
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 63.66262833333334
$ws.Cells.Item(6, 14).Value = 190.987885
$ws.Cells.Item(6, 15).Value = 0.3441717873742006
$ws.Cells.Item(6, 16).Value = 0.3441717873742006
$ws.Cells.Item(6, 17).Value = 794.7742471459818
$ws.Cells.Item(6, 18).Value = 7152.968224313836
$ws.Cells.Item(6, 19).Value = 0.3276670573984151
$ws.Cells.Item(6, 20).Value = 0.3276670573984151

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.484157
$ws.Cells.Item(7, 8).Value = 37.452471
$ws.Cells.Item(7, 9).Value = 0.952045081609665
$ws.Cells.Item(7, 10).Value = 0.9520450816096651
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.481835
$ws.Cells.Item(7, 14).Value = 55.445505
$ws.Cells.Item(7, 15).Value = 0.09991617299555507
$ws.Cells.Item(7, 16).Value = 0.09991617299555505
$ws.Cells.Item(7, 17).Value = 230.730129788095
$ws.Cells.Item(7, 18).Value = 2076.571168092855
$ws.Cells.Item(7, 19).Value = 0.09512470107367862
$ws.Cells.Item(7, 20).Value = 0.09512470107367862

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Sema3d"
$ws.Cells.Item(8, 3).Value = "Nrp1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5951003333333332
$ws.Cells.Item(8, 8).Value = 1.785301
$ws.Cells.Item(8, 9).Value = 0.04538250723811564
$ws.Cells.Item(8, 10).Value = 0.04538250723811565
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 102.8289443333334
$ws.Cells.Item(8, 14).Value = 308.486833
$ws.Cells.Item(8, 15).Value = 0.5559120396302444
$ws.Cells.Item(8, 16).Value = 0.5559120396302443
$ws.Cells.Item(8, 17).Value = 61.19353904908144
$ws.Cells.Item(8, 18).Value = 550.741851441733
$ws.Cells.Item(8, 19).Value = 0.0252286821622752
$ws.Cells.Item(8, 20).Value = 0.02522868216227519

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Sema3d"
$ws.Cells.Item(9, 3).Value = "Nrp1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5951003333333332
$ws.Cells.Item(9, 8).Value = 1.785301
$ws.Cells.Item(9, 9).Value = 0.04538250723811564
$ws.Cells.Item(9, 10).Value = 0.04538250723811565
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 63.66262833333334
$ws.Cells.Item(9, 14).Value = 190.987885
$ws.Cells.Item(9, 15).Value = 0.3441717873742006
$ws.Cells.Item(9, 16).Value = 0.3441717873742006
$ws.Cells.Item(9, 17).Value = 37.88565134204277
$ws.Cells.Item(9, 18).Value = 340.970862078385
$ws.Cells.Item(9, 19).Value = 0.01561937863166486
$ws.Cells.Item(9, 20).Value = 0.01561937863166486

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Sema3d"
$ws.Cells.Item(10, 3).Value = "Nrp1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5951003333333332
$ws.Cells.Item(10, 8).Value = 1.785301
$ws.Cells.Item(10, 9).Value = 0.04538250723811564
$ws.Cells.Item(10, 10).Value = 0.04538250723811565
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 18.481835
$ws.Cells.Item(10, 14).Value = 55.445505
$ws.Cells.Item(10, 15).Value = 0.09991617299555507
$ws.Cells.Item(10, 16).Value = 0.09991617299555505
$ws.Cells.Item(10, 17).Value = 10.99854616911166
$ws.Cells.Item(10, 18).Value = 98.98691552200498
$ws.Cells.Item(10, 19).Value = 0.004534446444175592
$ws.Cells.Item(10, 20).Value = 0.004534446444175592
